# Connector labels.xlsx - apply author's edit:
#   "Added second sensor outside the court. Given the physical condition
#    is on separate pin"
#
# Two logical changes:
#  1) Row 1's "^" (physical/ground condition symbol) and "9V" labels were
#     swapped between J1 and K1 so the condition symbol sits on its own
#     separate pin (matching the look of the other "^" cells in row 3),
#     and "9V" takes on the plain red label look that "^" used to have.
#  2) A second sensor's pin-label row (duplicate of row 1's A1:H1 block:
#     A0, A1, A2, A3, 5V, 5V, 5V, 5V) was added as row 5, leaving row 4
#     blank as a gap between the two sensors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Swap the "^" / "9V" labels on J1 / K1, carrying their formatting
#        with them so K1 ends up looking like the other "^" cells (A3:D3)
#        and J1 ends up looking like the old "9V" cell. ---
$ws.Range("K1").Copy($ws.Range("J1"))
$ws.Range("A3").Copy($ws.Range("K1"))

# --- 2) Add the second sensor's label row (copy of row 1's A1:H1) as
#        row 5, keeping row 4 empty as a spacer. ---
$ws.Range("A1:H1").Copy($ws.Range("A5:H5"))
